$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.431.28'
$ws.Range("E2").Value = '  -0.06%  '
$ws.Range("D3").Value = '1.848.65'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '240.84'
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("D6").Value = '0.6321'
$ws.Range("E6").Value = '  -3.57%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '3.225.29'
$ws.Range("E8").Value = '  +74.62%  '
$ws.Range("D9").Value = '0.07585'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").Value = '0.2970'
$ws.Range("E10").Value = '  -0.92%  '
$ws.Range("D11").Value = '24.58'
$ws.Range("E11").Value = '  +1.13%  '
$ws.Range("D12").Value = '0.07715'
$ws.Range("E12").Value = '  +1.10%  '
$ws.Range("D13").Value = '4.989'
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").Value = '0.6856'
$ws.Range("E14").Value = '  +0.16%  '
$ws.Range("D15").Value = '0.000009995'
$ws.Range("E15").Value = '  +4.92%  '
$ws.Range("D16").Value = '82.89'
$ws.Range("D17").Value = '6.187'
$ws.Range("D18").Value = '29.453.80'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '231.87'
$ws.Range("E19").Value = '  -2.38%  '
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").Value = '0.9997'
$ws.Range("E21").Value = '  +0.01%  '
$ws.Range("D22").Value = '7.577'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("D23").Value = '0.9998'
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '155.20'
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("D25").Value = '0.1389'
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("D26").Value = '8.438'
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("D28").Value = '1.470'
$ws.Range("E28").Value = '  -1.02%  '
$ws.Range("D29").Value = '0.05807'
$ws.Range("D30").Value = '1.259'
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").Value = '4.128'
$ws.Range("E31").Value = '  -0.17%  '
$ws.Range("D32").Value = '4.023'
$ws.Range("E32").Value = '  -1.27%  '
$ws.Range("D33").Value = '3.375.68'
$ws.Range("E33").Value = '  +68.73%  '
$ws.Range("D34").Value = '1.870'
$ws.Range("D35").Value = '1.158'
$ws.Range("E35").Value = '  -1.57%  '
$ws.Range("D36").Value = '0.7196'
$ws.Range("E36").Value = '  -0.29%  '
$ws.Range("E37").Value = '  -0.13%  '
$ws.Range("D38").Value = '1.248.72'
$ws.Range("E38").Value = '  +4.09%  '
$ws.Range("E39").Value = '  -0.21%  '
$ws.Range("D40").Value = '0.01808'
$ws.Range("E40").Value = '  +1.56%  '
$ws.Range("D41").Value = '0.9005'
$ws.Range("E41").Value = '  -1.17%  '
$ws.Range("D42").Value = '6.095'
$ws.Range("E42").Value = '  -2.17%  '
$ws.Range("D43").Value = '0.9993'
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("D44").Value = '101.37'
$ws.Range("E44").Value = '  -0.43%  '
$ws.Range("D46").Value = '7.323'
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").Value = '9.164'
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("D48").Value = '0.4019'
$ws.Range("D49").Value = '1.696'
$ws.Range("E49").Value = '  +2.17%  '
$ws.Range("D50").Value = '0.1126'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").Value = '0.05743'
$ws.Range("E51").Value = '  +0.25%  '
